$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Jianlei%Cao%NULL%0,                Xiaorong%Hu%NULL%0,                Wenlin%Cheng%NULL%0,                Lei%Yu%NULL%0,                Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                Qiang%Liu%liuqiang@irm-cams.ac.cn%0,                Qiang%Liu%liuqiang@irm-cams.ac.cn%0]"

$ws.Range("E4").Value = "[Chaolin%Huang%NULL%0,                Yeming%Wang%NULL%0,                Xingwang%Li%NULL%0,                Lili%Ren%NULL%0,                Jianping%Zhao%NULL%0,                Yi%Hu%NULL%0,                Li%Zhang%NULL%0,                Guohui%Fan%NULL%0,                Jiuyang%Xu%NULL%0,                Xiaoying%Gu%NULL%0,                Zhenshun%Cheng%NULL%0,                Ting%Yu%NULL%0,                Jiaan%Xia%NULL%0,                Yuan%Wei%NULL%0,                Wenjuan%Wu%NULL%0,                Xuelei%Xie%NULL%0,                Wen%Yin%NULL%0,                Hui%Li%NULL%0,                Min%Liu%NULL%0,                Yan%Xiao%NULL%0,                Hong%Gao%NULL%0,                Li%Guo%NULL%0,                Jungang%Xie%NULL%0,                Guangfa%Wang%NULL%0,                Rongmeng%Jiang%NULL%0,                Zhancheng%Gao%NULL%0,                Qi%Jin%NULL%0,                Jianwei%Wang%wangjw28@163.com%0,                Bin%Cao%caobin_ben@163.com%0]"

$ws.Range("E6").Value = "[Qiurong%Ruan%NULL%0,                Kun%Yang%NULL%0,                Kun%Yang%NULL%0,                Wenxia%Wang%NULL%0,                Wenxia%Wang%NULL%0,                Lingyu%Jiang%NULL%0,                Lingyu%Jiang%NULL%0,                Jianxin%Song%songsingsjx@sina.com%0,                Jianxin%Song%songsingsjx@sina.com%0]"

$ws.Range("E7").Value = "[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,               Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%3,               Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%3,               Subert-Salas%Lizandra%coreGivesNoEmail%3]"

$ws.Range("E8").Value = "[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,               Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%2,               Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%2,               Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%2,               Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%2]"

$ws.Range("E10").Value = "[Xiaobo%Yang%NULL%0,                Yuan%Yu%NULL%8,                Jiqian%Xu%NULL%8,                Huaqing%Shu%NULL%6,                Jia'an%Xia%NULL%0,                Hong%Liu%NULL%6,                Yongran%Wu%NULL%8,                Lu%Zhang%NULL%6,                Zhui%Yu%NULL%6,                Minghao%Fang%NULL%6,                Ting%Yu%NULL%0,                Yaxin%Wang%NULL%8,                Shangwen%Pan%NULL%6,                Xiaojing%Zou%NULL%6,                Shiying%Yuan%NULL%6,                You%Shang%NULL%8]"

$ws.Range("E11").Value = "[Fei%Zhou%NULL%0,                Ting%Yu%NULL%0,                Ronghui%Du%NULL%0,                Guohui%Fan%NULL%0,                Ying%Liu%NULL%0,                Zhibo%Liu%NULL%0,                Jie%Xiang%NULL%0,                Yeming%Wang%NULL%0,                Bin%Song%NULL%0,                Xiaoying%Gu%NULL%0,                Lulu%Guan%NULL%0,                Yuan%Wei%NULL%0,                Hui%Li%NULL%0,                Xudong%Wu%NULL%0,                Jiuyang%Xu%NULL%0,                Shengjin%Tu%NULL%0,                Yi%Zhang%NULL%0,                Hua%Chen%NULL%0,                Bin%Cao%NULL%0]"
